$d = $word.ActiveDocument

# Locate the "Scene design:" paragraph (the new content is inserted right
# after it, pushing the existing trailing empty paragraph further down
# without disturbing it).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    $t = $t.Replace([char]13, "").Replace([char]7, "")
    if ($t -eq "Scene design:") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Scene design:' paragraph"
}

# Build a range that covers the paragraph's text but stops just short of its
# paragraph mark, then collapse it to the end. Inserting XML exactly at the
# paragraph-mark boundary would consume/merge the following paragraph mark
# (observed via experimentation), so we back off by one character (the
# pilcrow) to keep insertion purely additive.
$tStart = $target.Range.Start
$tEnd = $target.Range.End - 1
$insertionRange = $d.Range($tStart, $tEnd)
$insertionRange.Collapse(0)  # wdCollapseEnd

$xmlFragment = '<w:p><w:r><w:t xml:space="preserve">The scene features a predominantly navy colour palette. The sky and the fog feature a similar colour, but not exactly the same in order to distinguish between the two. Considering the scene is a castle scene, a Victorian-style lamp has been used that illuminates a stack of wooden logs via a spotlight. </w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>3D techniques:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Fog has been implemented to create a hazy ambiance. A spotlight has been used to illuminate the logs in a way that contrasts the shadows in a pleasant way. The shadows are used to create a depth effect for the castle walls and </w:t></w:r><w:r><w:t>the logs. Normal mapping has been attempted but has not worked. Specular point lights illuminate the table. Several external OBJs with MTLs have been used to create a scene unique to the default ones used in previous levels.</w:t></w:r></w:p><w:p/><w:p><w:r><w:rPr><w:b/><w:bCs/><w:i/><w:iCs/><w:u w:val="single"/></w:rPr><w:t>Source code:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Multiple functions are implemented to render the shadows. The </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>OnRender</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t>) function has been moved around in the main.cpp in order to make the shadow functions work appropriately. Pre-initialisation was originally used, but moving the function proved to be more effective.</w:t></w:r></w:p><w:p/>'

$insertionRange.InsertXML($xmlFragment)
